# Script 1 - atualizacao automatica de dados
# Refreshes the quarterly 'Rendimento medio' series (Brasil / Nordeste / Sergipe)
# with the latest values and extends every region through 2025 Q3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$variavel = "Rendimento médio mensal real das pessoas de 14 anos ou mais de idade ocupadas na semana de referência com rendimento de trabalho, habitualmente recebido em todos os trabalhos"

# Each entry: Regiao, Trimestre, Valor
$rows = @(
    ,@("Brasil", "01/01/2017", 4639)
    ,@("Brasil", "01/04/2017", 4584)
    ,@("Brasil", "01/07/2017", 4599)
    ,@("Brasil", "01/10/2017", 4623)
    ,@("Brasil", "01/01/2018", 4501)
    ,@("Brasil", "01/04/2018", 4501)
    ,@("Brasil", "01/07/2018", 4494)
    ,@("Brasil", "01/10/2018", 4517)
    ,@("Brasil", "01/01/2019", 4357)
    ,@("Brasil", "01/04/2019", 4298)
    ,@("Brasil", "01/07/2019", 4307)
    ,@("Brasil", "01/10/2019", 4339)
    ,@("Brasil", "01/01/2020", 4210)
    ,@("Brasil", "01/04/2020", 4402)
    ,@("Brasil", "01/07/2020", 4452)
    ,@("Brasil", "01/10/2020", 4280)
    ,@("Brasil", "01/01/2021", 3843)
    ,@("Brasil", "01/04/2021", 3723)
    ,@("Brasil", "01/07/2021", 3579)
    ,@("Brasil", "01/10/2021", 3461)
    ,@("Brasil", "01/01/2022", 3330)
    ,@("Brasil", "01/04/2022", 3349)
    ,@("Brasil", "01/07/2022", 3467)
    ,@("Brasil", "01/10/2022", 3536)
    ,@("Brasil", "01/01/2023", 3405)
    ,@("Brasil", "01/04/2023", 3399)
    ,@("Brasil", "01/07/2023", 3462)
    ,@("Brasil", "01/10/2023", 3487)
    ,@("Brasil", "01/01/2024", 3371)
    ,@("Brasil", "01/04/2024", 3421)
    ,@("Brasil", "01/07/2024", 3415)
    ,@("Brasil", "01/10/2024", 3459)
    ,@("Brasil", "01/01/2025", 3367)
    ,@("Brasil", "01/04/2025", 3388)
    ,@("Brasil", "01/07/2025", 3406)
    ,@("Nordeste", "01/01/2017", 3144)
    ,@("Nordeste", "01/04/2017", 3128)
    ,@("Nordeste", "01/07/2017", 3072)
    ,@("Nordeste", "01/10/2017", 3171)
    ,@("Nordeste", "01/01/2018", 3083)
    ,@("Nordeste", "01/04/2018", 3065)
    ,@("Nordeste", "01/07/2018", 3077)
    ,@("Nordeste", "01/10/2018", 3090)
    ,@("Nordeste", "01/01/2019", 2976)
    ,@("Nordeste", "01/04/2019", 2951)
    ,@("Nordeste", "01/07/2019", 2921)
    ,@("Nordeste", "01/10/2019", 2953)
    ,@("Nordeste", "01/01/2020", 2865)
    ,@("Nordeste", "01/04/2020", 3024)
    ,@("Nordeste", "01/07/2020", 2930)
    ,@("Nordeste", "01/10/2020", 2847)
    ,@("Nordeste", "01/01/2021", 2533)
    ,@("Nordeste", "01/04/2021", 2504)
    ,@("Nordeste", "01/07/2021", 2411)
    ,@("Nordeste", "01/10/2021", 2333)
    ,@("Nordeste", "01/01/2022", 2214)
    ,@("Nordeste", "01/04/2022", 2212)
    ,@("Nordeste", "01/07/2022", 2294)
    ,@("Nordeste", "01/10/2022", 2326)
    ,@("Nordeste", "01/01/2023", 2283)
    ,@("Nordeste", "01/04/2023", 2277)
    ,@("Nordeste", "01/07/2023", 2287)
    ,@("Nordeste", "01/10/2023", 2314)
    ,@("Nordeste", "01/01/2024", 2249)
    ,@("Nordeste", "01/04/2024", 2352)
    ,@("Nordeste", "01/07/2024", 2313)
    ,@("Nordeste", "01/10/2024", 2361)
    ,@("Nordeste", "01/01/2025", 2318)
    ,@("Nordeste", "01/04/2025", 2323)
    ,@("Nordeste", "01/07/2025", 2338)
    ,@("Sergipe", "01/01/2017", 3643)
    ,@("Sergipe", "01/04/2017", 3483)
    ,@("Sergipe", "01/07/2017", 3325)
    ,@("Sergipe", "01/10/2017", 3230)
    ,@("Sergipe", "01/01/2018", 3254)
    ,@("Sergipe", "01/04/2018", 3103)
    ,@("Sergipe", "01/07/2018", 3244)
    ,@("Sergipe", "01/10/2018", 3164)
    ,@("Sergipe", "01/01/2019", 3008)
    ,@("Sergipe", "01/04/2019", 2950)
    ,@("Sergipe", "01/07/2019", 2918)
    ,@("Sergipe", "01/10/2019", 2836)
    ,@("Sergipe", "01/01/2020", 2897)
    ,@("Sergipe", "01/04/2020", 3085)
    ,@("Sergipe", "01/07/2020", 2960)
    ,@("Sergipe", "01/10/2020", 3126)
    ,@("Sergipe", "01/01/2021", 2561)
    ,@("Sergipe", "01/04/2021", 2759)
    ,@("Sergipe", "01/07/2021", 2586)
    ,@("Sergipe", "01/10/2021", 2513)
    ,@("Sergipe", "01/01/2022", 2257)
    ,@("Sergipe", "01/04/2022", 2265)
    ,@("Sergipe", "01/07/2022", 2368)
    ,@("Sergipe", "01/10/2022", 2425)
    ,@("Sergipe", "01/01/2023", 2353)
    ,@("Sergipe", "01/04/2023", 2404)
    ,@("Sergipe", "01/07/2023", 2320)
    ,@("Sergipe", "01/10/2023", 2288)
    ,@("Sergipe", "01/01/2024", 2249)
    ,@("Sergipe", "01/04/2024", 2372)
    ,@("Sergipe", "01/07/2024", 2399)
    ,@("Sergipe", "01/10/2024", 2594)
    ,@("Sergipe", "01/01/2025", 2480)
    ,@("Sergipe", "01/04/2025", 2436)
    ,@("Sergipe", "01/07/2025", 2713)
)

$headerRow = 1
$ws.Cells.Item($headerRow, 1).Value = "Região"
$ws.Cells.Item($headerRow, 2).Value = "Variável"
$ws.Cells.Item($headerRow, 3).Value = "Trimestre"
$ws.Cells.Item($headerRow, 4).Value = "Valor"

# Make sure the Trimestre column is forced to text, Excel would otherwise
# coerce the dd/mm/yyyy-looking strings into date serials.
$lastRow = $rows.Count + 1
$ws.Range("C2:C$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $region = $rows[$i][0]
    $quarter = $rows[$i][1]
    $value = $rows[$i][2]

    $ws.Cells.Item($r, 1).Value = $region
    $ws.Cells.Item($r, 2).Value = $variavel
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $quarter
    $ws.Cells.Item($r, 4).Value = $value
}
